$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Define the three new character styles that the commit introduces.
# ---------------------------------------------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2. Apply "GaNStyle" to every "Dates à utiliser..." run (4 occurrences).
# ---------------------------------------------------------------------------

$datesText = "Dates à utiliser pour la Campagne 2022 Constellation du Cygne: 10-19 août, 9-18 septembre, 8-17 octobre"
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = $datesText
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchWildcards = $false
while ($rng.Find.Execute()) {
    $rng.Style = "GaNStyle"
}

# ---------------------------------------------------------------------------
# 3. Apply "GaNParagraph" to the "Vous allez participer..." run.
# ---------------------------------------------------------------------------

$paragraphText = "Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation du Cygne dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne."
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = $paragraphText
$rng2.Find.Forward = $true
$rng2.Find.Wrap = 0
$rng2.Find.MatchWildcards = $false
while ($rng2.Find.Execute()) {
    $rng2.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------------
# 4. Apply "GaNLinks" to the "Les cartes figurant..." run.
# ---------------------------------------------------------------------------

$linksText = "Les cartes figurant dans ce document ont été établies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Text = $linksText
$rng3.Find.Forward = $true
$rng3.Find.Wrap = 0
$rng3.Find.MatchWildcards = $false
while ($rng3.Find.Execute()) {
    $rng3.Style = "GaNLinks"
}
